$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 00:22"

# --- Estados Unidos (row 4) updated totals ---
$ws.Range("B4").Value = 816336
$ws.Range("C4").Value = 23577
$ws.Range("D4").Value = 82676
$ws.Range("E4").Value = 688488
$ws.Range("G4").Value = 2658
$ws.Range("H4").Value = 45172

# --- Swap Finlandia / Colombia ordering (row 50 becomes Colombia w/ new data,
#     row 51 becomes Finlandia w/ what used to be row 50's data) ---
$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 4149
$ws.Range("C50").Value = 172
$ws.Range("D50").Value = 804
$ws.Range("E50").Value = 3149
$ws.Range("F50").Value = 98
$ws.Range("G50").Value = 7
$ws.Range("H50").Value = 196

$ws.Range("A51").Value = "Finlandia"
$ws.Range("B51").Value = 4014
$ws.Range("C51").Value = 146
$ws.Range("D51").Value = 2000
$ws.Range("E51").Value = 1873
$ws.Range("F51").Value = 63
$ws.Range("G51").Value = 43
$ws.Range("H51").Value = 141

# --- Uganda (row 158) updated figures ---
$ws.Range("C158").Value = 5
$ws.Range("E158").Value = 23
